$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a brand-new "Player Info" sheet as the first tab in the workbook.
# ---------------------------------------------------------------------------
$playerInfo = $wb.Worksheets.Add()
$playerInfo.Name = "Player Info"
$playerInfo.Move($wb.Worksheets.Item(1))

# Header row
$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

# Data row (ID is numeric-looking text, keep it as text like the rest of the workbook).
# Force text storage via NumberFormat, then drop the format again so the cell
# ends up with the plain/default style, exactly like every other data cell.
$playerInfo.Range("A2").NumberFormat = "@"
$playerInfo.Range("A2").Value = "6656"
$playerInfo.Range("A2").ClearFormats()
$playerInfo.Range("B2").Value = "Matthew Paul Kuhnemann"
$playerInfo.Range("C2").Value = "Left Handed"
$playerInfo.Range("D2").Value = "Left Arm Orthodox"

# Header styling - bold, centered, top-aligned, thin border all round (matches
# the look of the header rows already used on the other two sheets).
$headerRange = $playerInfo.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108   # xlCenter
$headerRange.VerticalAlignment = -4160     # xlTop
$headerRange.Borders.LineStyle = 1

$playerInfo.Range("A1").Select()

# ---------------------------------------------------------------------------
# 2. "ODI Batting" sheet - rename MATCH_CARD_LINK column to MATCH_CODE and
#    replace the full scorecard URL values with just the bare match code.
# ---------------------------------------------------------------------------
$odiBatting = $wb.Worksheets.Item("ODI Batting")
$odiBatting.Range("D1").Value = "MATCH_CODE"

$odiBatting.Range("D2").NumberFormat = "@"
$odiBatting.Range("D2").Value = "4597"
$odiBatting.Range("D2").ClearFormats()

$odiBatting.Range("D3").NumberFormat = "@"
$odiBatting.Range("D3").Value = "4600"
$odiBatting.Range("D3").ClearFormats()

$odiBatting.Range("D4").NumberFormat = "@"
$odiBatting.Range("D4").Value = "4601"
$odiBatting.Range("D4").ClearFormats()

$odiBatting.Range("D5").NumberFormat = "@"
$odiBatting.Range("D5").Value = "4603"
$odiBatting.Range("D5").ClearFormats()

# Re-apply the header style that ClearFormats could not have touched (D1 was
# never cleared) - left untouched on purpose, it already carries the bold
# header style from the original workbook.

# ---------------------------------------------------------------------------
# 3. "ODI Bowling" sheet - same MATCH_CARD_LINK -> MATCH_CODE change.
# ---------------------------------------------------------------------------
$odiBowling = $wb.Worksheets.Item("ODI Bowling")
$odiBowling.Range("B1").Value = "MATCH_CODE"

$odiBowling.Range("B2").NumberFormat = "@"
$odiBowling.Range("B2").Value = "4597"
$odiBowling.Range("B2").ClearFormats()

$odiBowling.Range("B3").NumberFormat = "@"
$odiBowling.Range("B3").Value = "4600"
$odiBowling.Range("B3").ClearFormats()

$odiBowling.Range("B4").NumberFormat = "@"
$odiBowling.Range("B4").Value = "4601"
$odiBowling.Range("B4").ClearFormats()

$odiBowling.Range("B5").NumberFormat = "@"
$odiBowling.Range("B5").Value = "4603"
$odiBowling.Range("B5").ClearFormats()
